$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich text runs) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "37"

$c9 = $ws.Range("C9")
$c9.Characters(27, 8).Text = "9/9/2024"
$c9.Characters(46, 8).Text = "9/15/2024"

# --- Crime statistics table updates (rows 15-28) ---
$ws.Range("G15").Value = 1
$ws.Range("N15").Value = -77.272727272727
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("I16").Value = 94
$ws.Range("J16").Value = 103
$ws.Range("K16").Value = -8.737864077669
$ws.Range("L16").Value = 23.684210526315
$ws.Range("M16").Value = -21.666666666666
$ws.Range("N16").Value = -82.877959927140
$ws.Range("C17").Value = 3
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 127
$ws.Range("J17").Value = 130
$ws.Range("K17").Value = -2.307692307692
$ws.Range("L17").Value = 14.414414414414
$ws.Range("M17").Value = 47.674418604651
$ws.Range("N17").Value = -39.234449760765
$ws.Range("C18").Value = 1
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -41.666666666666
$ws.Range("I18").Value = 40
$ws.Range("J18").Value = 72
$ws.Range("K18").Value = -44.444444444444
$ws.Range("L18").Value = -55.056179775280
$ws.Range("M18").Value = -28.571428571428
$ws.Range("N18").Value = -90.719257540603
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 21
$ws.Range("H19").Value = -34.375
$ws.Range("I19").Value = 239
$ws.Range("J19").Value = 283
$ws.Range("K19").Value = -15.547703180212
$ws.Range("L19").Value = -3.629032258064
$ws.Range("M19").Value = 27.127659574468
$ws.Range("N19").Value = -48.491379310344
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -66.666666666666
$ws.Range("J20").Value = 77
$ws.Range("K20").Value = -66.233766233766
$ws.Range("M20").Value = 52.941176470588
$ws.Range("N20").Value = -90.510948905109
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -29.411764705882
$ws.Range("G21").Value = 77
$ws.Range("H21").Value = -35.064935064935
$ws.Range("I21").Value = 532
$ws.Range("J21").Value = 672
$ws.Range("K21").Value = -20.833333333333
$ws.Range("L21").Value = -9.523809523809
$ws.Range("M21").Value = 10.602910602910
$ws.Range("N21").Value = -72.815533980582
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = -58.333333333333
$ws.Range("I23").Value = 94
$ws.Range("J23").Value = 111
$ws.Range("K23").Value = -15.315315315315
$ws.Range("L23").Value = 10.588235294117
$ws.Range("M23").Value = 56.666666666666
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = -42.857142857142
$ws.Range("F24").Value = 34
$ws.Range("G24").Value = 42
$ws.Range("H24").Value = -19.047619047619
$ws.Range("I24").Value = 290
$ws.Range("J24").Value = 363
$ws.Range("K24").Value = -20.110192837465
$ws.Range("L24").Value = -18.539325842696
$ws.Range("M24").Value = -36.681222707423
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = -27.272727272727
$ws.Range("I25").Value = 45
$ws.Range("J25").Value = 121
$ws.Range("K25").Value = -62.809917355371
$ws.Range("L25").Value = -62.809917355371
$ws.Range("C26").Value = 8
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 15
$ws.Range("H26").Value = 60
$ws.Range("I26").Value = 174
$ws.Range("J26").Value = 174
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 8.75
$ws.Range("M26").Value = -18.309859154929
$ws.Range("G27").Value = 1
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = -13.333333333333
